$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "prénom" (first name) column is repurposed into an "entreprise"
# (company) column, per "Mise a jour de l'envoi des mails".
$ws.Range("C3").Value = "entreprise  "
$ws.Range("C4").Value = "Thales"
$ws.Range("C5").Value = "Microsoft"

# Move the active selection, matching the author's last cursor position.
$ws.Range("K11").Select()
